$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column (D) receive numeric-looking text (e.g. "135.60",
# "64.301.50", "0.0000180"). Force Text format before assigning so Excel
# keeps the exact string instead of coercing it to a floating point number,
# matching the original inline-string cell content.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '64.301.50'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '3.502.46'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '584.63'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').Value = '135.60'
$ws.Range('E6').Value = '  +2.32%  '
$ws.Range('D7').Value = '3.503.69'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').Value = '  +0.28%  '
$ws.Range('D11').Value = '7.15'
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').Value = '0.375'
$ws.Range('E12').Value = '  -3.68%  '
$ws.Range('D13').Value = '4.096.08'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '0.0000180'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '3.503.38'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').Value = '25.74'
$ws.Range('E17').Value = '  -7.62%  '
$ws.Range('D18').Value = '64.284.48'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').Value = '9.78'
$ws.Range('E19').Value = '  -2.13%  '
$ws.Range('D20').Value = '13.84'
$ws.Range('E20').Value = '  -2.98%  '
$ws.Range('D21').Value = '5.57'
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range('D22').Value = '383.27'
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('E23').Value = '  -1.53%  '
$ws.Range('D24').Value = '3.637.27'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('D25').Value = '73.96'
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '5.71'
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('D28').Value = '0.0000115'
$ws.Range('E28').Value = '  +4.17%  '
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('D30').Value = '7.52'
$ws.Range('E30').Value = '  +1.39%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').Value = '8.29'
$ws.Range('E32').Value = '  +1.03%  '
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('D34').Value = '3.518.17'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D37').Value = '23.55'
$ws.Range('E37').Value = '  -1.84%  '
$ws.Range('D38').Value = '5.29'
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').Value = '1.55'
$ws.Range('E39').Value = '  -2.97%  '
$ws.Range('D40').Value = '6.84'
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('D41').Value = '163.96'
$ws.Range('E41').Value = '  -4.28%  '
$ws.Range('D42').Value = '0.0783'
$ws.Range('E42').Value = '  -3.41%  '
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').Value = '25.99'
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '41.96'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').Value = '1.21'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').Value = '4.40'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('D50').Value = '2.467.75'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '6.76'
$ws.Range('E51').Value = '  -2.13%  '

# Restore the default (Normal) style on those Price cells so the only
# recorded change is their text content, just like in the source edit.
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
